# Apply the "block order" reshuffle to the scenecat input file.
# The header labels get permuted across A1:F1, and the 0/1 indicator rows
# (rows 2-7) are updated so each row's single "1" marker moves to the
# column matching the new header ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column order ---
$headers = @("living_rooms_1", "living_rooms_2", "kitchens_1", "bedrooms_1", "bedrooms_2", "kitchens_2")
for ($col = 1; $col -le 6; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Data rows (rows 2-7): updated 0/1 indicator matrix ---
$data = @(
    @(0, 0, 0, 0, 1, 0),
    @(0, 0, 1, 0, 0, 0),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 1, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
